$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.7011554793401927
$ws.Range("C2").Value = 0.6931624611341978
$ws.Range("D2").Value = 0.7072593778714615
$ws.Range("E2").Value = 0.6815706391498696
$ws.Range("F2").Value = 0.7803321918335155

$ws.Range("B3").Value = 0.7304915890903152
$ws.Range("C3").Value = 0.7245803709759692
$ws.Range("D3").Value = 0.7372629361108668
$ws.Range("E3").Value = 0.7143244425120358
$ws.Range("F3").Value = 0.8117229101249219
